# RDU-48: rename the legacy "default"/"restricted" permission labels used
# throughout the excel2xml test workbook to the new standard names:
#   res-default / prop-default       -> open
#   res-restricted / prop-restricted -> restricted
#
# (All other cells / shared strings are unaffected - they just shift index
# in the shared-string table because these old labels are removed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$openCells = @(
    "G2",
    "O3","S3","W3",
    "O4","S4","W4","AA4",
    "O5","S5","W5","AA5",
    "O6","S6","W6","AA6",
    "O7","S7","W7","AA7",
    "O8",
    "G9",
    "O10",
    "O11",
    "G12","I12",
    "O13",
    "O14",
    "O15",
    "O16",
    "G17",
    "O18",
    "O19",
    "G25",
    "O26",
    "O27","S27","W27",
    "O28",
    "O29",
    "O30","S30","W30",
    "O31",
    "G32",
    "O33",
    "O34","S34"
)

$restrictedCells = @(
    "G20",
    "O21","S21",
    "O22",
    "O23",
    "O24",
    "S26"
)

foreach ($addr in $openCells) {
    $ws.Range($addr).Value = "open"
}

foreach ($addr in $restrictedCells) {
    $ws.Range($addr).Value = "restricted"
}

# The sheet's view was also scrolled/re-selected in the source edit.
$ws.Range("I7").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
